$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.352.72"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.242.68"
$ws.Range("E3").Value = "  +3.08%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.00"
$ws.Range("E5").Value = "  -1.89%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.68"
$ws.Range("E6").Value = "  -1.30%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.236.86"
$ws.Range("E8").Value = "  +2.95%  "
$ws.Range("E9").Value = "  -1.39%  "
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.34"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.38"
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.782.06"
$ws.Range("E15").Value = "  +3.20%  "
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.247.88"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "63.388.14"
$ws.Range("E18").Value = "  -0.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.78"
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "475.37"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.11"
$ws.Range("E21").Value = "  -4.20%  "
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.92"
$ws.Range("E23").Value = "  +3.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.04"
$ws.Range("E24").Value = "  -4.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.19"
$ws.Range("E25").Value = "  -0.85%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.44"
$ws.Range("E27").Value = "  +6.57%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.73"
$ws.Range("E28").Value = "  -1.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.11"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "27.52"
$ws.Range("E31").Value = "  -0.90%  "
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.53"
$ws.Range("E34").Value = "  -4.60%  "
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.90"
$ws.Range("E36").Value = "  -2.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.72"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0711"
$ws.Range("E38").Value = "  -3.89%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "421.73"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.39"
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.979.68"
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.75"
$ws.Range("E43").Value = "  -7.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.110"
$ws.Range("E44").Value = "  -7.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.267"
$ws.Range("E45").Value = "  +2.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.16"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("B48").Value = "ThetaToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.33"
$ws.Range("E48").Value = "  -2.97%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.86"
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "121.28"
$ws.Range("E51").Value = "  +0.77%  "
